$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.532.64"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "3.691.89"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "678.24"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.76"
$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.13"
$ws.Range("E10").Value = "  -0.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.440"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").Value = "4.312.26"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.52"
$ws.Range("E14").Value = "  +0.09%  "

$ws.Range("D15").Value = "3.713.00"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").Value = "69.458.59"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("E18").Value = "  +0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.81"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.81"
$ws.Range("E21").Value = "  -1.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.651"
$ws.Range("E22").Value = "  +0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "80.49"
$ws.Range("E23").Value = "  +1.19%  "

$ws.Range("D24").Value = "3.837.45"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("E27").Value = "  -1.08%  "

$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.71"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("E30").Value = "  -0.62%  "

$ws.Range("E31").Value = "  -0.42%  "

$ws.Range("E32").Value = "  -1.03%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.02"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").Value = "3.679.99"
$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("E36").Value = "  +1.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.49"
$ws.Range("E37").Value = "  +3.44%  "

$ws.Range("E38").Value = "  +1.01%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("E42").Value = "  -0.43%  "

$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.69"
$ws.Range("E45").Value = "  -2.52%  "

$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("E48").Value = "  -0.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "27.79"
$ws.Range("E49").Value = "  -2.50%  "

$ws.Range("E50").Value = "  -2.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.89"
$ws.Range("E51").Value = "  +0.86%  "
